$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tradeOwnerDetails")

$ws.Range("D2").Value = "Tester"
$ws.Range("E2").Value = "Tester Dad"
$ws.Range("F2").Value = "abc@xyz.com"
$ws.Range("C2").Value = 2222222222

$ws.Range("E10").Select()
